# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '28.310.79'
$ws.Range('E2').Value = '  +2.21%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.813.59'
$ws.Range('E3').Value = '  +3.37%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.000'
$ws.Range('E4').Value = '  -0.24%  '

# Row 5
Set-TextValue $ws.Range('D5') '325.52'
$ws.Range('E5').Value = '  +0.32%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.9988'
$ws.Range('E6').Value = '  -0.02%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.4371'
$ws.Range('E7').Value = '  +1.79%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.3663'
$ws.Range('E8').Value = '  +0.56%  '

# Row 9
Set-TextValue $ws.Range('D9') '44.67'
$ws.Range('E9').Value = '  -1.59%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.07660'

# Row 11
$ws.Range('E11').Value = '  +1.58%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.9993'
$ws.Range('E12').Value = '  -0.10%  '

# Row 13
Set-TextValue $ws.Range('D13') '21.98'
$ws.Range('E13').Value = '  +1.51%  '

# Row 14
Set-TextValue $ws.Range('D14') '6.298'
$ws.Range('E14').Value = '  +2.34%  '

# Row 15
Set-TextValue $ws.Range('D15') '7.484'
$ws.Range('E15').Value = '  +3.08%  '

# Row 16
Set-TextValue $ws.Range('D16') '1.827.64'
$ws.Range('E16').Value = '  +4.59%  '

# Row 17
Set-TextValue $ws.Range('D17') '95.07'
$ws.Range('E17').Value = '  +8.03%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.00001078'
$ws.Range('E18').Value = '  +0.90%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.06483'
$ws.Range('E19').Value = '  +4.47%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.9989'
$ws.Range('E20').Value = '  -0.03%  '

# Row 21
Set-TextValue $ws.Range('D21') '17.37'
$ws.Range('E21').Value = '  +1.46%  '

# Row 22
Set-TextValue $ws.Range('D22') '6.240'
$ws.Range('E22').Value = '  +1.25%  '

# Row 23
Set-TextValue $ws.Range('D23') '28.311.89'
$ws.Range('E23').Value = '  +2.12%  '

# Row 24
Set-TextValue $ws.Range('D24') '11.54'
$ws.Range('E24').Value = '  -1.28%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.106'
$ws.Range('E25').Value = '  -9.72%  '

# Row 26
Set-TextValue $ws.Range('D26') '161.34'
$ws.Range('E26').Value = '  +5.61%  '

# Row 27
Set-TextValue $ws.Range('D27') '20.72'
$ws.Range('E27').Value = '  +0.81%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.031.47'
$ws.Range('E28').Value = '  +4.29%  '

# Row 29
Set-TextValue $ws.Range('D29') '2.274'
$ws.Range('E29').Value = '  -4.13%  '

# Row 30
Set-TextValue $ws.Range('D30') '128.90'
$ws.Range('E30').Value = '  +1.16%  '

# Row 31
$ws.Range('E31').Value = '  -1.88%  '

# Row 32
Set-TextValue $ws.Range('D32') '6.003'
$ws.Range('E32').Value = '  +4.58%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.09144'
$ws.Range('E33').Value = '  -0.18%  '

# Row 34
Set-TextValue $ws.Range('D34') '3.581'
$ws.Range('E34').Value = '  -2.45%  '

# Row 35
Set-TextValue $ws.Range('D35') '12.92'
$ws.Range('E35').Value = '  +1.50%  '

# Row 36
Set-TextValue $ws.Range('D36') '0.02360'
$ws.Range('E36').Value = '  +1.97%  '

# Row 37
Set-TextValue $ws.Range('D37') '5.212'
$ws.Range('E37').Value = '  +1.76%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.2171'
$ws.Range('E38').Value = '  +0.63%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.6586'
$ws.Range('E39').Value = '  +1.45%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.06211'
$ws.Range('E40').Value = '  +1.68%  '

# Row 41
Set-TextValue $ws.Range('D41') '1.189'
$ws.Range('E41').Value = '  -0.63%  '

# Row 42
Set-TextValue $ws.Range('D42') '8.084'
$ws.Range('E42').Value = '  +1.28%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.423'
$ws.Range('E43').Value = '  -0.54%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.9985'
$ws.Range('E44').Value = '  -0.03%  '

# Row 45
Set-TextValue $ws.Range('D45') '13.77'
$ws.Range('E45').Value = '  -0.16%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.6103'
$ws.Range('E46').Value = '  +2.61%  '

# Row 47
$ws.Range('E47').Value = '  -0.46%  '

# Row 48
Set-TextValue $ws.Range('D48') '125.29'
$ws.Range('E48').Value = '  -0.75%  '

# Row 49
Set-TextValue $ws.Range('D49') '2.014'
$ws.Range('E49').Value = '  +2.01%  '

# Row 50
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D50') '1.154'
$ws.Range('E50').Value = '  +2.50%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.06990'
$ws.Range('E51').Value = '  +1.23%  '
